$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "admin"
$ws.Range("B1").Value = "manager"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

$ws.Range("B1:B2").Select()
